{"js": "// Change the acceptance-letter date from 25/02/16 to 29/02/16.\n// Word keeps an internal \"_GoBack\" bookmark marking the location of the\n// last edit; after changing the date, that bookmark moves from the\n// signature line (\"Sergio Avil\u00e9s\") to the date field itself (right\n// after the newly typed \"29\").\n\nconst doc = context.document;\nconst body = doc.body;\n\n// 1) Locate the day portion of the date (\"25\" inside \" 25/02/16\").\n//    It is unique in this document, so a plain search is safe.\nconst results = body.search(\"25\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const dayRange = results.items[0];\n\n  // 2) Replace \"25\" with \"29\"; insertText(\"Replace\") returns a range\n  //    over the freshly inserted text.\n  const newDayRange = dayRange.insertText(\"29\", \"Replace\");\n  await context.sync();\n\n  // 3) Move Word's \"_GoBack\" bookmark: delete it from its old spot\n  //    (after \"Sergio Avil\u00e9s\") and recreate it right after \"29\",\n  //    i.e. right before \"/02/16\" \u2014 matching where the edit happened.\n  doc.deleteBookmark(\"_GoBack\");\n  const afterDay = newDayRange.getRange(\"After\");\n  afterDay.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Change the acceptance-letter date from 25/02/16 to 29/02/16.\n# Word keeps an internal \"_GoBack\" bookmark marking the location of the\n# last edit; after changing the date, that bookmark moves from the\n# signature line (\"Sergio Avil\u00e9s\") to the date field itself (right\n# after the newly typed \"29\").\n\n$d = $word.ActiveDocument\n\n# 1) Locate \" 25/02/16\" (unique in this document).\n$rng = $d.Content\n$rng.Find.Text = \"25/02/16\"\n$rng.Find.MatchCase = $true\n$rng.Find.Forward = $true\n$found = $rng.Find.Execute()\n\nif ($found) {\n    $s = $rng.Start\n    $e = $rng.End\n\n    # 2) Pin the already-existing run boundary right before the leading\n    #    space (i.e. between \"Fecha:\" and \" 25/02/16\") with a scratch\n    #    bookmark so the in-place text edit below only rewrites the day\n    #    digits and doesn't get folded back into the \"Fecha:\" run.\n    $d.Range($s - 1, $s - 1).Bookmarks.Add(\"_TempPin\") | Out-Null\n\n    # 3) Move Word's \"_GoBack\" bookmark to sit right after the new day\n    #    value, i.e. between \"29\" and \"/02/16\" -- adding a bookmark under\n    #    a name that already exists relocates it from wherever it was\n    #    before (after \"Sergio Avil\u00e9s\").\n    $d.Range($s + 2, $s + 2).Bookmarks.Add(\"_GoBack\") | Out-Null\n\n    # 4) Change the day from 25 to 29.\n    $d.Range($s, $s + 2).Text = \"29\"\n\n    # 5) Drop the temporary pin now that the edit is done.\n    $d.Bookmarks.Item(\"_TempPin\").Delete()\n}\n"}
